$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "A03"
$ws.Range("B4").Value = "Shiguelose"
$ws.Range("A5").Value = "A04"
$ws.Range("B5").Value = "Outras infecções intestinais bacterianas"
$ws.Range("A6").Value = "A09"
$ws.Range("B6").Value = "Diarréia e gastroenterite de origem infecciosa presumível"
$ws.Range("A7").Value = "A90"
$ws.Range("B7").Value = "Dengue [dengue clássico]"
$ws.Range("A8").Value = "B50"
$ws.Range("B8").Value = "Malária por Plasmodium falciparum"
$ws.Range("A9").Value = "B54"
$ws.Range("B9").Value = "Malária não especificada"
$ws.Range("A10").Value = "B34.2"
$ws.Range("B10").Value = "Infecção por coronavírus de localização não especificada"
$ws.Range("A11").Value = "J00"
$ws.Range("B11").Value = "Nasofaringite aguda [resfriado comum]"
$ws.Range("A12").Value = "J01"
$ws.Range("B12").Value = "Sinusite aguda"
$ws.Range("A13").Value = "J02"
$ws.Range("B13").Value = "Faringite aguda"
$ws.Range("A14").Value = "J03"
$ws.Range("B14").Value = "Amigdalite aguda"
$ws.Range("A15").Value = "J06"
$ws.Range("B15").Value = "Infecções agudas das vias aéreas superiores"
$ws.Range("A16").Value = "J10"
$ws.Range("B16").Value = "Influenza devida a vírus da influenza identificado"
$ws.Range("A17").Value = "J11"
$ws.Range("B17").Value = "Influenza devida a vírus não identificado"
$ws.Range("A18").Value = "J12"
$ws.Range("B18").Value = "Pneumonia viral"
$ws.Range("A19").Value = "J18"
$ws.Range("B19").Value = "Pneumonia por microorganismo não especificado"
$ws.Range("A20").Value = "J18.9"
$ws.Range("B20").Value = "Pneumonia não especificada"
$ws.Range("A21").Value = "J20"
$ws.Range("B21").Value = "Bronquite aguda"
$ws.Range("A22").Value = "J44"
$ws.Range("B22").Value = "Outras doenças pulmonares obstrutivas crônicas"
$ws.Range("A23").Value = "J45"
$ws.Range("B23").Value = "Asma"
$ws.Range("A24").Value = "I10"
$ws.Range("B24").Value = "Hipertensão essencial (primária)"
$ws.Range("A25").Value = "I20"
$ws.Range("B25").Value = "Angina pectoris"
$ws.Range("A26").Value = "I21"
$ws.Range("B26").Value = "Infarto agudo do miocárdio"
$ws.Range("A27").Value = "I50"
$ws.Range("B27").Value = "Insuficiência cardíaca"
$ws.Range("A28").Value = "E10"
$ws.Range("B28").Value = "Diabetes mellitus insulino-dependente"
$ws.Range("A29").Value = "E11"
$ws.Range("B29").Value = "Diabetes mellitus não-insulino-dependente"
$ws.Range("A30").Value = "E66"
$ws.Range("B30").Value = "Obesidade"
$ws.Range("A31").Value = "F32"
$ws.Range("B31").Value = "Episódios depressivos"
$ws.Range("A32").Value = "F41"
$ws.Range("B32").Value = "Outros transtornos ansiosos"
$ws.Range("A33").Value = "M54"
$ws.Range("B33").Value = "Dorsalgia"
$ws.Range("A34").Value = "M54.5"
$ws.Range("B34").Value = "Dor lombar baixa"
$ws.Range("A35").Value = "R50"
$ws.Range("B35").Value = "Febre de outra origem e de origem desconhecida"
$ws.Range("A36").Value = "R51"
$ws.Range("B36").Value = "Cefaléia"
$ws.Range("A37").Value = "R52"
$ws.Range("B37").Value = "Dor não classificada em outra parte"
$ws.Range("A38").Value = "S82"
$ws.Range("B38").Value = "Fratura da perna, incluindo tornozelo"
$ws.Range("A39").Value = "Z00.0"
$ws.Range("B39").Value = "Exame geral e investigação de pessoas sem queixas ou diagnóstico relatado"
$ws.Range("A40").Value = "Z76"
$ws.Range("B40").Value = "Pessoas em contato com serviços de saúde em outras circunstâncias"
